# Diving_log_BenthFun.xlsx edit
# Commit message: LOW TN2 + O2 PROFILES T3 + TN2 + LIGHT PROFILE T3 + TN2 (MISSING ELOW)
#
# Work happens on the "Corrected" worksheet (second tab, tabSelected):
#  1. Rows 311-326 (Tn2_t1_ELOW_* tiles): date in column A moves from
#     2023-09-10 (45187) to 2023-09-11 (45188).
#  2. Rows 327-342 (Tn2_t1_LOW_* tiles + O2 profiles T3/Tn2): date in column A
#     moves from 2023-09-08 (45185) to 2023-09-10 (45187), and the
#     previously-missing Start/End/Bottom time (D/E/F) and Temperature/Light
#     (I/J) measurements are filled in.
#  3. The active selection moves to L337 with the view scrolled further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Corrected")

# --- 1. Rows 311-326: bump the dive date by one day (45187 -> 45188) -------
foreach ($r in 311..326) {
    $ws.Cells.Item($r, 1).Value2 = 45188
}

# --- 2. Rows 327-342: bump the dive date (45185 -> 45187) and backfill the
#        previously-missing D/E/F (times) and I/J (temperature / light) data
$rowsData = @(
    @{Row=327; D=0.45624999999999999; E=0.50277777777777777; F=0.50277777777777777; I=26.378;              J=-71.2},
    @{Row=328; D=0.45694444444444443; E=0.50624999999999998; F=0.50624999999999998; I=26.513000000000002;  J=-70.2},
    @{Row=329; D=0.45763888888888887; E=0.50972222222222219; F=0.50972222222222219; I=26.77;               J=-71.8},
    @{Row=330; D=0.45833333333333331; E=0.5131944444444444;  F=0.5131944444444444;  I=26.411999999999999;  J=-55.5},
    @{Row=331; D=0.46597222222222223; E=0.49791666666666662; F=0.49791666666666662; I=25.733000000000001;  J=-59.3},
    @{Row=332; D=0.46736111111111112; E=0.4993055555555555;  F=0.4993055555555555;  I=25.838000000000001;  J=-56.5},
    @{Row=333; D=0.4680555555555555;  E=0.5;                 F=0.5;                 I=25.841999999999999;  J=-58.9},
    @{Row=334; D=0.45833333333333331; E=0.5131944444444444;  F=0.5131944444444444;  I=26.411999999999999;  J=-55.5},
    @{Row=335; D=0.52708333333333335; E=0.57291666666666663; F=0.57291666666666663; I=26.817;              J=-75.900000000000006},
    @{Row=336; D=0.52777777777777779; E=0.57500000000000007; F=0.57500000000000007; I=26.91;               J=-73.2},
    @{Row=337; D=0.52847222222222223; E=0.57847222222222217; F=0.57847222222222217; I=27.146000000000001;  J=-77.8},
    @{Row=338; D=0.52916666666666667; E=0.58194444444444449; F=0.58194444444444449; I=27.056000000000001;  J=-62},
    @{Row=339; D=0.53611111111111109; E=0.56805555555555554; F=0.56805555555555554; I=26.405000000000001;  J=-60.9},
    @{Row=340; D=0.53680555555555554; E=0.56944444444444442; F=0.56944444444444442; I=26.393000000000001;  J=-62.5},
    @{Row=341; D=0.53749999999999998; E=0.56944444444444442; F=0.56944444444444442; I=26.509;              J=-64.400000000000006},
    @{Row=342; D=0.52916666666666667; E=0.58194444444444449; F=0.58194444444444449; I=27.056000000000001;  J=-62}
)

foreach ($rd in $rowsData) {
    $r = $rd.Row
    $ws.Cells.Item($r, 1).Value2 = 45187
    $ws.Cells.Item($r, 4).Value2 = $rd.D
    $ws.Cells.Item($r, 5).Value2 = $rd.E
    $ws.Cells.Item($r, 6).Value2 = $rd.F
    $ws.Cells.Item($r, 9).Value2 = $rd.I
    $ws.Cells.Item($r, 10).Value2 = $rd.J
}

# --- 3. Update the active sheet view / selection ----------------------------
$ws.Activate()
$ws.Range("L337").Select()
